$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Number total of files" (column H) for rows whose metric was recomputed
$ws.Cells.Item(5, 8).Value = 2955
$ws.Cells.Item(7, 8).Value = 60208
$ws.Cells.Item(8, 8).Value = 3443
$ws.Cells.Item(11, 8).Value = 2292
$ws.Cells.Item(13, 8).Value = 2535
$ws.Cells.Item(14, 8).Value = 4213
$ws.Cells.Item(15, 8).Value = 4538
$ws.Cells.Item(18, 8).Value = 1376
$ws.Cells.Item(19, 8).Value = 1516
$ws.Cells.Item(20, 8).Value = 2195
$ws.Cells.Item(21, 8).Value = 1084
$ws.Cells.Item(22, 8).Value = 3268
$ws.Cells.Item(24, 8).Value = 91115
$ws.Cells.Item(26, 8).Value = 3214
$ws.Cells.Item(28, 8).Value = 2947
$ws.Cells.Item(29, 8).Value = 4667
$ws.Cells.Item(31, 8).Value = 5212
$ws.Cells.Item(32, 8).Value = 2410
$ws.Cells.Item(33, 8).Value = 6138
$ws.Cells.Item(34, 8).Value = 3833
$ws.Cells.Item(35, 8).Value = 3685
$ws.Cells.Item(36, 8).Value = 12770
$ws.Cells.Item(37, 8).Value = 583
$ws.Cells.Item(38, 8).Value = 4498
$ws.Cells.Item(40, 8).Value = 6514
$ws.Cells.Item(42, 8).Value = 1878
$ws.Cells.Item(43, 8).Value = 3093
$ws.Cells.Item(44, 8).Value = 5337
$ws.Cells.Item(45, 8).Value = 2913
$ws.Cells.Item(46, 8).Value = 5481
$ws.Cells.Item(47, 8).Value = 6235
$ws.Cells.Item(48, 8).Value = 21305
$ws.Cells.Item(49, 8).Value = 5014
$ws.Cells.Item(52, 8).Value = 22487
$ws.Cells.Item(53, 8).Value = 4505
$ws.Cells.Item(54, 8).Value = 15596
$ws.Cells.Item(55, 8).Value = 14562
$ws.Cells.Item(56, 8).Value = 10718
$ws.Cells.Item(57, 8).Value = 33499
$ws.Cells.Item(58, 8).Value = 11057
$ws.Cells.Item(63, 8).Value = 724
$ws.Cells.Item(64, 8).Value = 155056
$ws.Cells.Item(65, 8).Value = 9560
$ws.Cells.Item(66, 8).Value = 2226
$ws.Cells.Item(68, 8).Value = 8354
$ws.Cells.Item(69, 8).Value = 8675
$ws.Cells.Item(71, 8).Value = 3622
$ws.Cells.Item(78, 8).Value = 7043
$ws.Cells.Item(79, 8).Value = 3886
$ws.Cells.Item(80, 8).Value = 3499
$ws.Cells.Item(82, 8).Value = 12287
$ws.Cells.Item(84, 8).Value = 16612
$ws.Cells.Item(86, 8).Value = 2785
$ws.Cells.Item(87, 8).Value = 6983
$ws.Cells.Item(91, 8).Value = 1149
$ws.Cells.Item(92, 8).Value = 1174
$ws.Cells.Item(94, 8).Value = 978
$ws.Cells.Item(95, 8).Value = 24069
$ws.Cells.Item(97, 8).Value = 5642
$ws.Cells.Item(103, 8).Value = 9023
$ws.Cells.Item(104, 8).Value = 3121
$ws.Cells.Item(105, 8).Value = 1085
$ws.Cells.Item(109, 8).Value = 538
$ws.Cells.Item(111, 8).Value = 4290
$ws.Cells.Item(112, 8).Value = 3277
$ws.Cells.Item(113, 8).Value = 1138
$ws.Cells.Item(115, 8).Value = 1081
$ws.Cells.Item(116, 8).Value = 31704
$ws.Cells.Item(119, 8).Value = 797
$ws.Cells.Item(120, 8).Value = 1287
$ws.Cells.Item(121, 8).Value = 618
$ws.Cells.Item(122, 8).Value = 7941
$ws.Cells.Item(123, 8).Value = 2837
$ws.Cells.Item(124, 8).Value = 4920
$ws.Cells.Item(125, 8).Value = 421
$ws.Cells.Item(126, 8).Value = 1274
$ws.Cells.Item(127, 8).Value = 977
$ws.Cells.Item(128, 8).Value = 4062
$ws.Cells.Item(129, 8).Value = 15358
$ws.Cells.Item(130, 8).Value = 21197
$ws.Cells.Item(131, 8).Value = 1368
$ws.Cells.Item(132, 8).Value = 1118
$ws.Cells.Item(133, 8).Value = 608
$ws.Cells.Item(134, 8).Value = 18421
$ws.Cells.Item(135, 8).Value = 17397
$ws.Cells.Item(136, 8).Value = 4705
$ws.Cells.Item(137, 8).Value = 1071
$ws.Cells.Item(141, 8).Value = 314
$ws.Cells.Item(143, 8).Value = 931
$ws.Cells.Item(144, 8).Value = 3891
$ws.Cells.Item(145, 8).Value = 2423
$ws.Cells.Item(151, 8).Value = 89712
$ws.Cells.Item(154, 8).Value = 1329
$ws.Cells.Item(155, 8).Value = 1143
$ws.Cells.Item(156, 8).Value = 10104
$ws.Cells.Item(157, 8).Value = 2097
$ws.Cells.Item(158, 8).Value = 815
$ws.Cells.Item(160, 8).Value = 12452
$ws.Cells.Item(161, 8).Value = 1546
$ws.Cells.Item(165, 8).Value = 4453
$ws.Cells.Item(166, 8).Value = 432
$ws.Cells.Item(167, 8).Value = 2344
$ws.Cells.Item(168, 8).Value = 1629
$ws.Cells.Item(169, 8).Value = 2986
$ws.Cells.Item(173, 8).Value = 7799
$ws.Cells.Item(175, 8).Value = 3347
$ws.Cells.Item(177, 8).Value = 148
$ws.Cells.Item(178, 8).Value = 1737
$ws.Cells.Item(182, 8).Value = 2585
$ws.Cells.Item(187, 8).Value = 12442
$ws.Cells.Item(188, 8).Value = 1342
$ws.Cells.Item(189, 8).Value = 31764
$ws.Cells.Item(191, 8).Value = 5079
$ws.Cells.Item(195, 8).Value = 1574
$ws.Cells.Item(197, 8).Value = 963
$ws.Cells.Item(200, 8).Value = 1906
$ws.Cells.Item(201, 8).Value = 1141
$ws.Cells.Item(203, 8).Value = 9808
$ws.Cells.Item(204, 8).Value = 4975
$ws.Cells.Item(205, 8).Value = 2052
$ws.Cells.Item(209, 8).Value = 3705
$ws.Cells.Item(210, 8).Value = 3655
$ws.Cells.Item(216, 8).Value = 4205
$ws.Cells.Item(217, 8).Value = 1466
$ws.Cells.Item(218, 8).Value = 3316
$ws.Cells.Item(220, 8).Value = 8658
$ws.Cells.Item(223, 8).Value = 36993
$ws.Cells.Item(224, 8).Value = 1670
$ws.Cells.Item(225, 8).Value = 3942
$ws.Cells.Item(226, 8).Value = 3903
$ws.Cells.Item(228, 8).Value = 1468
$ws.Cells.Item(229, 8).Value = 7963
$ws.Cells.Item(230, 8).Value = 3083
$ws.Cells.Item(233, 8).Value = 1103
$ws.Cells.Item(234, 8).Value = 4051
$ws.Cells.Item(236, 8).Value = 3598
$ws.Cells.Item(238, 8).Value = 1668
$ws.Cells.Item(240, 8).Value = 20870
$ws.Cells.Item(242, 8).Value = 1357
$ws.Cells.Item(243, 8).Value = 4793
$ws.Cells.Item(244, 8).Value = 1582
$ws.Cells.Item(245, 8).Value = 3404
$ws.Cells.Item(246, 8).Value = 1339
$ws.Cells.Item(247, 8).Value = 2589
$ws.Cells.Item(251, 8).Value = 50432
$ws.Cells.Item(253, 8).Value = 488
$ws.Cells.Item(254, 8).Value = 2217
$ws.Cells.Item(256, 8).Value = 4068
$ws.Cells.Item(257, 8).Value = 437
$ws.Cells.Item(260, 8).Value = 1759
$ws.Cells.Item(261, 8).Value = 16619
$ws.Cells.Item(263, 8).Value = 632
$ws.Cells.Item(264, 8).Value = 7656
$ws.Cells.Item(266, 8).Value = 1029
$ws.Cells.Item(267, 8).Value = 10561
$ws.Cells.Item(268, 8).Value = 3126
$ws.Cells.Item(271, 8).Value = 4714
$ws.Cells.Item(272, 8).Value = 4909
$ws.Cells.Item(274, 8).Value = 11954
$ws.Cells.Item(275, 8).Value = 7996
$ws.Cells.Item(277, 8).Value = 1009
$ws.Cells.Item(280, 8).Value = 1642
$ws.Cells.Item(282, 8).Value = 24653
$ws.Cells.Item(284, 8).Value = 4464
$ws.Cells.Item(288, 8).Value = 1528
$ws.Cells.Item(289, 8).Value = 696
$ws.Cells.Item(290, 8).Value = 314
$ws.Cells.Item(292, 8).Value = 1712
$ws.Cells.Item(293, 8).Value = 2924
$ws.Cells.Item(297, 8).Value = 447
$ws.Cells.Item(299, 8).Value = 8265
$ws.Cells.Item(300, 8).Value = 2716
$ws.Cells.Item(302, 8).Value = 1262
$ws.Cells.Item(303, 8).Value = 5794
$ws.Cells.Item(304, 8).Value = 1199
$ws.Cells.Item(305, 8).Value = 7757
$ws.Cells.Item(307, 8).Value = 2145
$ws.Cells.Item(310, 8).Value = 11222
$ws.Cells.Item(311, 8).Value = 4129
$ws.Cells.Item(313, 8).Value = 1295
$ws.Cells.Item(314, 8).Value = 1560
$ws.Cells.Item(315, 8).Value = 3701
$ws.Cells.Item(316, 8).Value = 849
$ws.Cells.Item(317, 8).Value = 816
$ws.Cells.Item(318, 8).Value = 1087
$ws.Cells.Item(321, 8).Value = 741
$ws.Cells.Item(322, 8).Value = 992
$ws.Cells.Item(323, 8).Value = 2468
$ws.Cells.Item(324, 8).Value = 4677
$ws.Cells.Item(326, 8).Value = 10196
$ws.Cells.Item(327, 8).Value = 3008
$ws.Cells.Item(328, 8).Value = 3003
$ws.Cells.Item(329, 8).Value = 2678
$ws.Cells.Item(330, 8).Value = 1138
$ws.Cells.Item(332, 8).Value = 2868
$ws.Cells.Item(333, 8).Value = 3137
$ws.Cells.Item(334, 8).Value = 7410
$ws.Cells.Item(335, 8).Value = 1339
$ws.Cells.Item(336, 8).Value = 5014
$ws.Cells.Item(337, 8).Value = 8172
$ws.Cells.Item(338, 8).Value = 5327
$ws.Cells.Item(339, 8).Value = 1517
$ws.Cells.Item(340, 8).Value = 14724
$ws.Cells.Item(341, 8).Value = 1444
$ws.Cells.Item(345, 8).Value = 1551
$ws.Cells.Item(346, 8).Value = 2462
$ws.Cells.Item(349, 8).Value = 289
$ws.Cells.Item(351, 8).Value = 22073
$ws.Cells.Item(356, 8).Value = 13495
$ws.Cells.Item(357, 8).Value = 2921

# Append new row 363 (liferay/liferay-portal) to the results table
$ws.Cells.Item(363, 1).Value = "liferay"
$ws.Cells.Item(363, 2).Value = "liferay-portal"
$ws.Cells.Item(363, 3).Value = "Enterprise Resource Planning"
$ws.Cells.Item(363, 4).Value = "38220df7f113ab58af84b3297d81db6bb40d3320"
$ws.Cells.Item(363, 5).Value = ""
$ws.Cells.Item(363, 6).Value = "2024-10-12 00:30:11+00:00"
$ws.Cells.Item(363, 7).Value = $true
$ws.Cells.Item(363, 8).Value = 109041
$ws.Cells.Item(363, 9).Value = 0
$ws.Cells.Item(363, 10).Value = 0
$ws.Cells.Item(363, 11).Value = 32
$ws.Cells.Item(363, 12).Value = 215
$ws.Cells.Item(363, 13).Value = 155
$ws.Cells.Item(363, 14).Value = 1
$ws.Cells.Item(363, 15).Value = 1332
